$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainFont($rng) {
    $rng.Font.Bold = $false
    $rng.Font.Italic = $false
    $rng.Font.Strikethrough = $false
    $rng.Font.Underline = -4142
}

# --- Phase 1: write new translated values in the same order they were
# originally authored (this determines sharedStrings.xml insertion order) ---
$ws.Range("C18").Value = "ED1　True Ending"
$ws.Range("C6").Value = "\n<Lily>What's keeping them?`nIt's been such a long time."
$ws.Range("C8").Value = "\n<Lily>Ah- Right there♥`nThat hits the spot."
$ws.Range("C11").Value = "\n<Lily>Huh!?`nNo way!!`nAre you kidding me!?"
$ws.Range("C14").Value = "\n<Lily>I want to chase after them right away, but...`nGrumble...`nI need to put you two back to normal first..."
$ws.Range("C17").Value = "\n<Lily>Yes...`nI'll chase them to the ends of the earth...`nThen I'll drain them into dust...!!"
$ws.Range("C7").Value = "\n<Lime>Yep-."
$ws.Range("C9").Value = "\n<Lime>Right here-?"
$ws.Range("C15").Value = "\n<Lime>I don't want to keep losing either-.`nIt's so frustrating."
$ws.Range("C10").Value = "\n<Shina>Oiii!! There's no one here, nya!!`nThe dungeon's empty too!!`nAnd I found some rope dangling from the balcony, nya!!"
$ws.Range("C12").Value = "\n<Shina>I ain't kidding, nyan...`nAko is missing too.`nThat bastard... She betrayed us, nyan!"
$ws.Range("C16").Value = "\n<Shina>How dare they beat me and then escape, nya...`nI've never been so humiliated, nyan.`nThey're dead next time."
$ws.Range("C13").Value = "\n<Lily>\n[1]...!!`nDamn you...!!`nI'm gonna remember this!!"

# --- Phase 2: formatting (font, wrap, row heights, column widths) ---
$ws.Columns.Item(1).ColumnWidth = 41.736979166666664
$ws.Columns.Item(2).ColumnWidth = 45.307291666666664
$ws.Columns.Item(3).ColumnWidth = 70.59244791666666

# Row 6
Set-PlainFont($ws.Range("C6"))
$ws.Range("C6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 30

# Row 7
Set-PlainFont($ws.Range("C7"))

# Row 8
Set-PlainFont($ws.Range("C8"))
$ws.Range("C8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 30

# Row 9
Set-PlainFont($ws.Range("C9"))

# Row 10
Set-PlainFont($ws.Range("C10"))
$ws.Range("C10").WrapText = $true
Set-PlainFont($ws.Range("B10"))
$ws.Range("B10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 60

# Row 11
Set-PlainFont($ws.Range("C11"))
$ws.Range("C11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 45

# Row 12
Set-PlainFont($ws.Range("C12"))
$ws.Range("C12").WrapText = $true
Set-PlainFont($ws.Range("B12"))
$ws.Range("B12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 45

# Row 13
Set-PlainFont($ws.Range("C13"))
$ws.Range("C13").WrapText = $true
Set-PlainFont($ws.Range("B13"))
$ws.Range("B13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 45

# Row 14
Set-PlainFont($ws.Range("C14"))
$ws.Range("C14").WrapText = $true
Set-PlainFont($ws.Range("B14"))
$ws.Range("B14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 60

# Row 15
Set-PlainFont($ws.Range("C15"))
$ws.Range("C15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 30

# Row 16
Set-PlainFont($ws.Range("C16"))
$ws.Range("C16").WrapText = $true
Set-PlainFont($ws.Range("B16"))
$ws.Range("B16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 45

# Row 17
Set-PlainFont($ws.Range("C17"))
$ws.Range("C17").WrapText = $true
Set-PlainFont($ws.Range("B17"))
$ws.Range("B17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 45

# Row 18
Set-PlainFont($ws.Range("C18"))

# --- Phase 3: selection / view state ---
$ws.Range("C17").Select()
$excel.ActiveWindow.DisplayGridlines = $true